$d = $word.ActiveDocument

# --- Change 1: merge the two runs of the intro sentence into one run ---
# Re-assigning Range.Text rewrites the covered runs as a single run (using the
# first run's formatting), collapsing the original two runs
# ("...should" + " be included (maximum 2 pages).") into one. Unlike
# Find.Execute's Replacement text, a direct Range.Text assignment is not
# mangled by smart-quote autocorrect, so the straight apostrophe survives.
$introPara = $d.Paragraphs.Item(1)
$introRange = $d.Range($introPara.Range.Start, $introPara.Range.End - 1)
$introRange.Text = "Any strong individually held views on some aspect of the project, that the group doesn't agree on, should be included (maximum 2 pages)."

# --- Change 2: add Edward's new quote paragraph after the "Edward" heading ---
$q1 = [char]0x201C   # "
$q2 = [char]0x2019   # '
$q3 = [char]0x201D   # "

# Locate the "Edward" paragraph, then the (currently empty) paragraph right after it.
$edwardPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Edward") {
        $edwardPara = $d.Paragraphs.Item($i)
        break
    }
}
$targetIndex = $edwardPara.Index + 1
$targetPara = $d.Paragraphs.Item($targetIndex)

$r = $targetPara.Range
$r.Collapse(1)   # wdCollapseStart

$quoteBody = "works well together and I feel we are able to be very productive. Unfortunately I don" + $q2 + "t feel that we put enough time into planning the inner workings of the system in the early stages of the project and as such met slight problems when it came to the integration of the project. We were able to produce a system that works well for a project of this scale, but if it were to grow substantially then a slight refactor would probably be needed - I don" + $q2 + "t think this would be a huge task, but it would probably require us to step back and rethink some of our strategies. On a whole our group managed to get the work done, however I do feel that I put more time into the project than the other members of the group; this is however probably mainly to do with the fact I am somewhat of a perfectionist. I feel that this project has been a great learning experience, as I have learned what I am good at, and what I need to improve at." + $q3

$fullText = $q1 + "Our group " + $quoteBody
$r.InsertAfter($fullText)

$pStart = $targetPara.Range.Start
$paraText = $targetPara.Range.Text

# Split point right after the opening quote mark (between "“" and "Our group ").
$splitAfterQuote = $pStart + 1
$tmpName = "zzTmpSplit"
$d.Bookmarks.Add($tmpName, $d.Range($splitAfterQuote, $splitAfterQuote)) | Out-Null
$d.Bookmarks.Item($tmpName).Delete()

# Split point right after "Our group " (before "works well together...") -- this is
# also where the _GoBack bookmark belongs; adding it both splits the run and
# (re)places the bookmark here, automatically removing it from its old location
# in Brian's paragraph since bookmark names must be unique in the document.
$idx = $paraText.IndexOf("works well together")
$splitBeforeBody = $pStart + $idx
$d.Bookmarks.Add("_GoBack", $d.Range($splitBeforeBody, $splitBeforeBody)) | Out-Null
